$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 79499
$ws.Range("B3").Value = 91828
$ws.Range("B4").Value = 80349
$ws.Range("B5").Value = 80348
$ws.Range("B6").Value = 57884
